# "Added validation on generate button"
#
# Both sheets ("rice" and "wheat") get a new "serial number" column
# inserted immediately before the existing "From" column (their current
# column B), pushing "From".."Commodity"/"Values" one column to the right.
#
# On the "wheat" sheet the trailing "Values" header is additionally
# dropped (so the row only grows from B1:G1 to B1:G1, but shifted), and a
# sample data row is appended right under the header.

$wb = $excel.ActiveWorkbook

$xlShiftToRight = [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight
$xlCenter       = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$xlTop          = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignTop
$xlContinuous   = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$xlThin         = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

function Format-HeaderCell($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlTop
    $cell.Borders.LineStyle = $xlContinuous
    $cell.Borders.Weight = $xlThin
}

# ---------------------------------------------------------------------
# Sheet "rice": insert a new "Sr. No" column before "From".
# ---------------------------------------------------------------------
$wsRice = $wb.Worksheets.Item("rice")

$wsRice.Range("B1").Insert($xlShiftToRight)
$wsRice.Range("B1").Value = "Sr. No"
Format-HeaderCell $wsRice.Range("B1")

# ---------------------------------------------------------------------
# Sheet "wheat": insert a new "sr" column before "From", drop the old
# trailing "Values" header, and add one sample data row.
# ---------------------------------------------------------------------
$wsWheat = $wb.Worksheets.Item("wheat")

# "Values" is not kept around (it would otherwise get shifted into H1).
$wsWheat.Range("G1").Clear()

$wsWheat.Range("B1").Insert($xlShiftToRight)
$wsWheat.Range("B1").Value = "sr"
Format-HeaderCell $wsWheat.Range("B1")

$wsWheat.Range("A2").Value = 0
Format-HeaderCell $wsWheat.Range("A2")

$wsWheat.Range("C2").Value = "ENB"
$wsWheat.Range("D2").Value = "Haryana"
$wsWheat.Range("E2").Value = "BBU"
$wsWheat.Range("F2").Value = "Bihar"
$wsWheat.Range("G2").Value = "Wheat"
